$d = $word.ActiveDocument

# Insert a new, empty paragraph at the very start of the document body.
# The new paragraph carries only paragraph-mark formatting (bold) and no
# run/text, matching:
#   <w:p>
#     <w:pPr>
#       <w:rPr>
#         <w:b w:val="1"/>
#         <w:bCs w:val="1"/>
#       </w:rPr>
#     </w:pPr>
#   </w:p>
$r = $d.Range(0, 0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr></w:pPr></w:p>'
$r.InsertXML($xml)
